$d = $word.ActiveDocument

# Locate the "For 6 people:" paragraph robustly via Find.
$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute("For 6 people:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $para = $findRange.Paragraphs(1)
    $full = $para.Range

    # Apply bold + italic (regular and complex-script variants) to the whole
    # paragraph, including its end-of-paragraph mark, so the paragraph mark
    # ends up carrying the new "bold italic" run formatting.
    $full.Font.Bold = $true
    $full.Font.Italic = $true
    $full.Font.BoldBi = $true
    $full.Font.ItalicBi = $true

    # Re-create the visible run text (minus the trailing colon) from scratch
    # so the bold/italic formatting we just applied doesn't linger on the
    # visible text run -- only the paragraph mark should keep it.
    $textOnly = $d.Range($full.Start, $full.End - 1)
    $textOnly.Delete()
    $ins = $d.Range($full.Start, $full.Start)
    $ins.InsertBefore("For 6 people")
}
